$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1481481481481481
$ws.Range("C2").Value = 0.6666666666666666
$ws.Range("J2").Value = 0.02314814814814815
$ws.Range("P2").Value = 0.07407407407407407
$ws.Range("S2").Value = 0.08796296296296297
$ws.Range("C3").Value = 0.03333333333333333
$ws.Range("J3").Value = 0.02666666666666667
$ws.Range("O3").Value = 0.006666666666666667
$ws.Range("P3").Value = 0.7866666666666666
$ws.Range("S3").Value = 0.1466666666666667
$ws.Range("J4").Value = 0.07894736842105263
$ws.Range("P4").Value = 0.5789473684210527
$ws.Range("S4").Value = 0.3421052631578947
$ws.Range("B6").Value = 0.03720930232558139
$ws.Range("F6").Value = 0.03255813953488372
$ws.Range("J6").Value = 0.3162790697674419
$ws.Range("O6").Value = 0.02325581395348837
$ws.Range("Q6").Value = 0.1767441860465116
$ws.Range("R6").Value = 0.06976744186046512
$ws.Range("S6").Value = 0.3441860465116279
$ws.Range("B7").Value = 0.08695652173913043
$ws.Range("D7").Value = 0.02717391304347826
$ws.Range("E7").Value = 0.005434782608695652
$ws.Range("F7").Value = 0.03804347826086957
$ws.Range("J7").Value = 0.1467391304347826
$ws.Range("O7").Value = 0.02717391304347826
$ws.Range("Q7").Value = 0.2065217391304348
$ws.Range("R7").Value = 0.09239130434782608
$ws.Range("S7").Value = 0.3695652173913043
$ws.Range("B8").Value = 0.08351648351648351
$ws.Range("D8").Value = 0.01978021978021978
$ws.Range("E8").Value = 0.002197802197802198
$ws.Range("F8").Value = 0.07472527472527472
$ws.Range("J8").Value = 0.08791208791208792
$ws.Range("O8").Value = 0.01098901098901099
$ws.Range("Q8").Value = 0.1824175824175824
$ws.Range("R8").Value = 0.0989010989010989
$ws.Range("S8").Value = 0.4395604395604396
$ws.Range("B9").Value = 0.0963855421686747
$ws.Range("D9").Value = 0.01204819277108434
$ws.Range("E9").Value = 0.006024096385542169
$ws.Range("F9").Value = 0.03614457831325301
$ws.Range("J9").Value = 0.108433734939759
$ws.Range("O9").Value = 0.01204819277108434
$ws.Range("Q9").Value = 0.1807228915662651
$ws.Range("R9").Value = 0.1385542168674699
$ws.Range("S9").Value = 0.4096385542168675
$ws.Range("B10").Value = 0.09325044404973357
$ws.Range("D10").Value = 0.02042628774422735
$ws.Range("F10").Value = 0.07282415630550622
$ws.Range("J10").Value = 0.1021314387211368
$ws.Range("O10").Value = 0.01776198934280639
$ws.Range("Q10").Value = 0.2255772646536412
$ws.Range("R10").Value = 0.1039076376554174
$ws.Range("S10").Value = 0.3641207815275311
$ws.Range("G11").Value = 0.1564885496183206
$ws.Range("J11").Value = 0.07251908396946564
$ws.Range("K11").Value = 0.1908396946564886
$ws.Range("L11").Value = 0.5572519083969466
$ws.Range("S11").Value = 0.02290076335877863
$ws.Range("G12").Value = 0.7763157894736842
$ws.Range("J12").Value = 0.1578947368421053
$ws.Range("K12").Value = 0.006578947368421052
$ws.Range("L12").Value = 0.03289473684210526
$ws.Range("S12").Value = 0.02631578947368421
$ws.Range("G13").Value = 0.7631578947368421
$ws.Range("J13").Value = 0.1842105263157895
$ws.Range("S13").Value = 0.05263157894736842
$ws.Range("F15").Value = 0.02577319587628866
$ws.Range("H15").Value = 0.1443298969072165
$ws.Range("I15").Value = 0.06701030927835051
$ws.Range("J15").Value = 0.3814432989690721
$ws.Range("K15").Value = 0.06185567010309279
$ws.Range("O15").Value = 0.05154639175257732
$ws.Range("S15").Value = 0.2680412371134021
$ws.Range("F16").Value = 0.006666666666666667
$ws.Range("H16").Value = 0.2
$ws.Range("I16").Value = 0.06
$ws.Range("J16").Value = 0.38
$ws.Range("K16").Value = 0.09333333333333334
$ws.Range("M16").Value = 0.01333333333333333
$ws.Range("O16").Value = 0.08
$ws.Range("S16").Value = 0.1666666666666667
$ws.Range("F17").Value = 0.02921348314606742
$ws.Range("H17").Value = 0.1685393258426966
$ws.Range("I17").Value = 0.0696629213483146
$ws.Range("J17").Value = 0.4179775280898876
$ws.Range("K17").Value = 0.08764044943820225
$ws.Range("M17").Value = 0.01573033707865169
$ws.Range("N17").Value = 0.002247191011235955
$ws.Range("O17").Value = 0.07415730337078652
$ws.Range("S17").Value = 0.1348314606741573
$ws.Range("F18").Value = 0.02314814814814815
$ws.Range("H18").Value = 0.2037037037037037
$ws.Range("I18").Value = 0.1064814814814815
$ws.Range("J18").Value = 0.375
$ws.Range("K18").Value = 0.1203703703703704
$ws.Range("M18").Value = 0.01388888888888889
$ws.Range("O18").Value = 0.03240740740740741
$ws.Range("S18").Value = 0.125
$ws.Range("F19").Value = 0.02339688041594454
$ws.Range("H19").Value = 0.2435008665511265
$ws.Range("I19").Value = 0.0779896013864818
$ws.Range("J19").Value = 0.3552859618717504
$ws.Range("K19").Value = 0.1013864818024263
$ws.Range("M19").Value = 0.02426343154246101
$ws.Range("N19").Value = 0.0008665511265164644
$ws.Range("O19").Value = 0.06239168110918544
$ws.Range("S19").Value = 0.1109185441941074
